$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current header row (row 75), pushing the
# header (was row 75) and the two existing data rows (were rows 76-77) down
# by one to rows 76-78.
$ws.Rows.Item(75).Insert()

# --- New data row 79: sg_rr_36_025, prominence 0.015 run ---
$ws.Range("A79").Value = "sg_rr_36_025 2023-12-13 16-41-08.csv"
$ws.Range("B79").Value = 0.01
$ws.Range("C79").Value = 1000
$ws.Range("D79").Value = 5001
$ws.Range("E79").Value = 1530
$ws.Range("F79").Value = 1570
$ws.Range("G79").Value = 0.015
$ws.Range("H79").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I79").Value = 2.5
$ws.Range("U79").Value = "looks like prominence was probably too low as seems visually to find peaks in noise, going to try increasing it again without recording rest of data."
$ws.Range("V79").Value = "looks like prominence was probably too low as seems visually to find peaks in noise, going to try increasing it again without recording rest of data."

# --- New data row 80: sg_rr_36_025, prominence 0.02 run ---
$ws.Range("A80").Value = "sg_rr_36_025 2023-12-13 16-41-08.csv"
$ws.Range("B80").Value = 0.01
$ws.Range("C80").Value = 1000
$ws.Range("D80").Value = 5001
$ws.Range("E80").Value = 1530
$ws.Range("F80").Value = 1570
$ws.Range("G80").Value = 0.02
$ws.Range("H80").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I80").Value = 2.5
$ws.Range("U80").Value = "looks like prominence was probably too low as seems visually to find peaks in noise, going to try increasing it again without recording rest of data."
$ws.Range("V80").Value = "looks like prominence was probably too low as seems visually to find peaks in noise, going to try increasing it again without recording rest of data."

# --- New data row 81: sg_rr_36_025, prominence 0.021 run (full results) ---
$ws.Range("A81").Value = "sg_rr_36_025 2023-12-13 16-41-08.csv"
$ws.Range("B81").Value = 0.01
$ws.Range("C81").Value = 1000
$ws.Range("D81").Value = 5001
$ws.Range("E81").Value = 1530
$ws.Range("F81").Value = 1570
$ws.Range("G81").Value = 0.021
$ws.Range("H81").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I81").Value = 2.5
$ws.Range("J81").Value = 2.72785714285714
$ws.Range("K81").Value = 0.0133909339278385
$ws.Range("L81").Value = "yes"
$ws.Range("M81").Value = 0.156375206779422
$ws.Range("N81").Value = 0.00722041202165073
$ws.Range("O81").Value = 10163.8309034193
$ws.Range("P81").Value = 389.618530121367
$ws.Range("Q81").Value = 810153460.869798
$ws.Range("R81").Value = 93277498.0423375
$ws.Range("S81").Value = 36
$ws.Range("T81").Value = 0.1

# Update the sheet view: scroll so row 68 is at the top and select the new
# final row A81 (matching where data entry left off).
$ws.Application.ActiveWindow.ScrollRow = 68
$ws.Range("A81").Select()
